$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $find"
    }
}

# 1. Ativação date change
# NOTE: a plain Find/Replace here would cause the engine to merge the edited
# run with the immediately-following run (since both are plain, unformatted
# runs and the final run in the paragraph lacks a trailing <w:br/>). To keep
# the original run/break structure intact (as in the target), rebuild the
# whole "Créditos..." paragraph via InsertXML with the exact desired markup.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Créditos-aula: 4")) {
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListNumber"/></w:pPr><w:r><w:t>Créditos-aula: 4</w:t><w:br/></w:r><w:r><w:t>Créditos-trabalho: 0</w:t><w:br/></w:r><w:r><w:t>Carga horária: 60 h</w:t><w:br/></w:r><w:r><w:t>Semestre ideal: 6</w:t><w:br/></w:r><w:r><w:t>Ativação: 01/01/2021</w:t><w:br/></w:r><w:r><w:t>Departamento: Engenharia Química</w:t></w:r></w:p>'
        [void]$p.Range.InsertXML($xml)
        break
    }
}

# 2. Objetivos paragraph (PT) - remove "aos alunos de Engenharia de Produção"
Replace-Text "Apresentar aos alunos de Engenharia de Produção os conceitos básicos da Ciência Econômica, capacitando-os a compreender os principais conceitos micro e macroeconômicos e a interpretar o discurso e a prática da economia, orientados pelo seu próprio senso crítico." "Apresentar  os conceitos básicos da Ciência Econômica, capacitando-os a compreender os principais conceitos micro e macroeconômicos e a interpretar o discurso e a prática da economia, orientados pelo seu próprio senso crítico."

# 3. Objetivos paragraph (EN)
Replace-Text "Introduce students of production engineering the basic concepts of Economic Science, enabling them to understand the main micro and macroeconomic concepts and to interpret the discourse and practice of economics, guided by their own critical sense." "Introduce the students of Production Engineering to the basic concepts of Economic Science, enabling students to understand the main micro and macroeconomic concepts and to interpret the discourse and practice of economics, guided by their own critical sense."

# 4. Remove the "5840671 - Francisco José Moreira Chaves" paragraph entirely
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "5840671 - Francisco José Moreira Chaves") {
        $p.Range.Delete()
        break
    }
}

# 5. "Programa resumido" PT
Replace-Text "1  Microecomonia. 2 - Macroeconomia. 3 - Análise da Economia Brasileira" "A. Microeconomia. B. Macroeconomia. C. Desenvolvimento Econômico. D. Economia Internacional. E. Economia Brasileira"

# 6. "Programa resumido" EN
Replace-Text "1 - Microeconomics. 2 - Macroeconomics. 3 - Analysis of the Brazilian Economy" "A. Microeconomy. B. Macroeconomy. C. Economic Development. D. International Economy. E. Brazilian Economy"

# 7. "Programa" PT (long)
Replace-Text "MICROECONOMIA1. Introdução aos conceitos de Economia e fundamentos da análise microeconômica.2. Teoria do consumidor e da demanda.3. Teoria da firma e da oferta. 4. Custos e formação de preços. 5. Estruturas de Mercado6. Comportamento estratégico e concorrência.7. Tecnologia como fator de produção.8. Sustentabilidade: recursos, custos e indicadores ambientais.MACROECONOMIA1. Fundamentos da análise macroeconômica.2. Contabilidade nacional.3. Equilíbrios clássicos e keynesiano.4. Sistema monetário.5. Política fiscal.6. Economia mundial e comércio internacional.7. Fundamentos da regressão como ferramenta para quantificar relações econômicas.8. Setor público.ECONOMIA BRASILEIRA1. A experiência histórica da industrialização brasileira.2. A internacionalização da economia brasileira.3. Teoria dos ciclos e realidade brasileira.4. Os ciclos econômicos do Brasil ao longo de sua história recente." "A. MICROECONOMIA: 1. Introdução aos conceitos de Economia e fundamentos da análise microeconômica. 2. Teoria do consumidor e da demanda. 3. Teoria da firma e da oferta. 4. Custos e formação de preços. 5. Estruturas de Mercado 6. Comportamento estratégico e concorrência. 7. Tecnologia como fator de produção. 8. Sustentabilidade: recursos, custos e indicadores ambientais. B. MACROECONOMIA: 1. Fundamentos da análise macroeconômica. 2. Contabilidade nacional. 3. Equilíbrios clássicos e keynesiano. 4. Sistema monetário. 5. Política fiscal. 6. Economia mundial e comércio internacional. 7. Fundamentos da regressão como ferramenta para quantificar relações econômicas. 8. Setor público. C. DESENVOLVIMENTO ECONÔMICO: 1. Fatores de Crescimento. 2. Fontes de Desenvolvimento. 3. Financiamento do Desenvolvimento Econômico. 4. Um modelo de Crescimento Econômico. 5. O Processo de internacionalização e globalização.D. ECONOMIA INTERNACIONAL: 1. Fundamentos do Comércio Internacional. 2. Determinação das Taxas de Câmbio. 3. Políticas Externas. 4. Fatores determinantes do comportamento das importações e exportações.E. ECONOMIA BRASILEIRA: 1. A experiência histórica da industrialização brasileira. 2. A internacionalização da economia brasileira. 3. Teoria dos ciclos e realidade brasileira. 4. Os ciclos econômicos do Brasil ao longo de sua história recente."

# 8. "Programa" EN (long)
# NOTE: this text contains a plain straight apostrophe ("Brazil's"). A
# Find/Replace would run it through the engine's smart-quote autocorrect and
# turn it into a curly apostrophe (’), which would not match the target.
# InsertXML bypasses autocorrect and writes the literal text as given.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("MICROECONOMICS1. Introduction")) {
        $apos = [char]39
        $xmlStart = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>A. MICROECONOMY: 1. Introduction to the concepts of Economics and fundamentals of microeconomic analysis. 2. Consumer and demand theory. 3. Firm and offer theory. 4. Costs and price formation. 5. Market Structures 6. Strategic behavior and competition. 7. Technology as a factor of production. 8. Sustainability: resources, costs and environmental indicators.B. MACROECONOMY: 1. Fundamentals of macroeconomic analysis. 2. National accounting. 3. Classical and Keynesian balances. 4. Monetary system. 5. Fiscal policy. 6. World economy and international trade. 7. Fundamentals of regression as a tool to quantify economic relationships. 8. Public sector.C. ECONOMIC DEVELOPMENT: 1. Growth factors. 2. Sources of Development. 3. Financing Economic Development. 4. A model of economic growth. 5. The internationalization and globalization process.D. INTERNATIONAL ECONOMY: 1. Fundamentals of International Trade. 2. Determination of Exchange Rates. 3. External policies. 4. Factors determining the behavior of imports and exports.E. BRAZILIAN ECONOMY: 1. The historical experience of Brazilian industrialization. 2. The internationalization of the Brazilian economy. 3. Cycle theory and Brazilian reality. 4. Brazil'
        $xmlEnd = 's economic cycles throughout its recent history.</w:t></w:r></w:p>'
        $xml = $xmlStart + $apos + $xmlEnd
        [void]$p.Range.InsertXML($xml)
        break
    }
}

# 9. Método (Avaliação)
Replace-Text "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras." "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

# 10. Critério (Avaliação)
Replace-Text "MF = (0,40*P1 + 0,40*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários." "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas."

# 11. Norma de recuperação (Avaliação)
Replace-Text "NF = (MF + PR)/2, onde PR é uma prova de recuperação" "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# 12. Bibliografia - replace entire paragraph content, removing all <w:br/> separators
# Using InsertXML (instead of setting Range.Text) produces a clean <w:t> run
# without the engine tacking on an xml:space="preserve" attribute that isn't
# present in the target markup.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("BACHA , Edmar. Introdução")) {
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>VASCONCELLOS, M. A. S.; GARCIA, M. E. Fundamentos de Economia. 6 ed. São Paulo: Saraiva, 2018.GREMAUD, A. P. Introdução à Economia. São Paulo: Atlas, 2017.ROSSETTI, J. P. Introdução à Economia - Livro Texto. São Paulo: Atlas, 2016.VASCONCELLOS, M. A. S. ECONOMIA: Micro e Macro. São Paulo: Atlas, 2015.ALBERGONI, L. INTRODUÇÃO À ECONOMIA: Aplicações no Cotidiano. São Paulo: Atlas, 2015.GREMAUD, A. P.; VASCONCELLOS, M. A. S.; TONETO JÚNIOR, R. Economia Brasileira Contemporânea. 8 ed. São Paulo: Atlas, 2017.MÉNARD, C.; SAES, M. S. M.; SILVA, V. L. S.; RAYNAUD, E. Economia das Organizações: Formas Plurais e Desafios. São Paulo: Atlas, 2014.BACHA et Al. Estado da Economia Mundial - Desafios e Respostas - Seminário em Homenagem a Pedro Malan. São Paulo: LTC, 2015.BACHA , Edmar. Introdução à Macroeconomia: Uma perspectiva brasileira. Rio de Janeiro: Campus,1987.BEGG, D.; DORNBUSCH, R.; FISCHER, S. Introdução A Economia. Rio de Janeiro: Campus, 2003. FURTADO, C. Formação econômica do Brasil. São Paulo: Companhia Editora Nacional, 2003.GRAMAUD, A. P. et alli. Manual de economia. São Paulo. Saraiva. 2004.GRAMAUD, A. P. et alli. Economia Brasileira Contemporânea. 6.ed. São Paulo. Atlas, 2006.HUNT, E. K.; SHERMAN, H. J. História do Pensamento Econômico. Petrópolis : Vozes, 1997.MANKIW, N.G. Introdução à economia. São Paulo: Thomson Learning, 2006.SAMUELSON, P. Introdução à Economia. New York: Mc Graw-Hill Book Company.</w:t></w:r></w:p>'
        [void]$p.Range.InsertXML($xml)
        break
    }
}

Write-Host "Done"
